# Fixes Weekly Tutorial case typo.
# On slide 11, the paragraph "There are Weekly tutorials that are required..."
# has the word "tutorials" (run 4) in lower-case even though it is formatted
# bold alongside "Weekly". Capitalize it to "Tutorials" to match.

$p = $ppt.ActivePresentation
$s = $p.Slides.Item(11)
$shape = $s.Shapes.Item(2)
$tf = $shape.TextFrame
$tr = $tf.TextRange

# Paragraph 2 is "There are Weekly tutorials that are required to be
# completed by students for a 2% grade for each tutorial. ..."
$para = $tr.Paragraphs(2, 1)

# Locate the bold run whose text is exactly "tutorials" and fix its case.
$fixed = $false
for ($i = 1; $i -le $para.Runs().Count; $i++) {
    $run = $para.Runs($i, 1)
    if ($run.Text -ceq "tutorials") {
        $run.Text = "Tutorials"
        $fixed = $true
        break
    }
}

if (-not $fixed) {
    Write-Host "WARNING: target run 'tutorials' not found - no changes made"
}
